$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.307.30'
$ws.Range('E2').Value = '  -4.99%  '

$ws.Range('D3').Value = '3.138.43'
$ws.Range('E3').Value = '  -5.02%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.20%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '515.96'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -7.66%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.94'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.51%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('D8').Value = '3.138.04'
$ws.Range('E8').Value = '  -5.09%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.444'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -6.80%  '

$ws.Range('E10').Value = '  -8.44%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.107'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -9.94%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.380'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -6.74%  '

$ws.Range('D13').Value = '3.668.94'
$ws.Range('E13').Value = '  -4.83%  '

$ws.Range('E14').Value = '  -2.16%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.11'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -7.62%  '

$ws.Range('D16').Value = '3.133.83'
$ws.Range('E16').Value = '  -4.98%  '

$ws.Range('D17').Value = '57.339.66'
$ws.Range('E17').Value = '  -4.79%  '

$ws.Range('E18').Value = '  -10.90%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.70'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -6.84%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.79'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -10.65%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.88'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -8.02%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '341.39'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -9.17%  '

$ws.Range('E23').Value = '  -0.18%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.34'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -7.43%  '

$ws.Range('E25').Value = '  -8.15%  '

$ws.Range('D26').Value = '3.260.87'

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.13%  '

$ws.Range('E28').Value = '  -6.15%  '

$ws.Range('D29').Value = '0.0₃0924'
$ws.Range('E29').Value = '  -10.91%  '

$ws.Range('E30').Value = '  -0.04%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.64'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -8.23%  '

$ws.Range('E32').Value = '  -9.82%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '21.45'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.88%  '

$ws.Range('E34').Value = '  -11.09%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.17'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.19%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.78'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -8.40%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '157.31'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.56%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.14'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -8.90%  '

$ws.Range('E39').Value = '  -9.88%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '25.55'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.08%  '

$ws.Range('D41').Value = '3.161.31'
$ws.Range('E41').Value = '  -4.89%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0679'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -8.54%  '

$ws.Range('E43').Value = '  -3.83%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.689'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -8.23%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.05'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.99%  '

$ws.Range('E46').Value = '  +0.18%  '

$ws.Range('E47').Value = '  -8.49%  '

$ws.Range('E48').Value = '  -9.13%  '

$ws.Range('D49').Value = '2.229.00'
$ws.Range('E49').Value = '  -5.27%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.09'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.61%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.81'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.53%  '
